$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "27.178.77"
Set-TextValue $ws "E2" "  -2.45%  "

# Row 3
Set-TextValue $ws "D3" "1.866.53"
Set-TextValue $ws "E3" "  -2.10%  "

# Row 4
Set-TextValue $ws "E4" "  -0.26%  "

# Row 5
Set-TextValue $ws "D5" "307.28"
Set-TextValue $ws "E5" "  -1.98%  "

# Row 6
Set-TextValue $ws "D6" "1.001"
Set-TextValue $ws "E6" "  -0.19%  "

# Row 7
Set-TextValue $ws "D7" "0.5158"
Set-TextValue $ws "E7" "  +2.87%  "

# Row 8
Set-TextValue $ws "D8" "0.3751"

# Row 9
Set-TextValue $ws "D9" "0.07172"
Set-TextValue $ws "E9" "  -1.51%  "

# Row 10
Set-TextValue $ws "B10" "Solana"
Set-TextValue $ws "C10" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws "D10" "20.67"
Set-TextValue $ws "E10" "  -0.79%  "

# Row 11
Set-TextValue $ws "B11" "Polygon"
Set-TextValue $ws "C11" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D11" "0.8845"
Set-TextValue $ws "E11" "  -2.53%  "

# Row 12
Set-TextValue $ws "B12" "TRON"
Set-TextValue $ws "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D12" "0.07566"
Set-TextValue $ws "E12" "  -1.33%  "

# Row 13
Set-TextValue $ws "B13" "WrappedEther"
Set-TextValue $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D13" "1.863.61"
Set-TextValue $ws "E13" "  -2.38%  "

# Row 14
Set-TextValue $ws "D14" "5.336"
Set-TextValue $ws "E14" "  -2.66%  "

# Row 15
Set-TextValue $ws "E15" "  -2.66%  "

# Row 16
Set-TextValue $ws "D16" "1.001"
Set-TextValue $ws "E16" "  -0.22%  "

# Row 17
Set-TextValue $ws "D17" "0.000008566"
Set-TextValue $ws "E17" "  -1.74%  "

# Row 18
Set-TextValue $ws "D18" "14.14"
Set-TextValue $ws "E18" "  -2.59%  "

# Row 19
Set-TextValue $ws "D19" "1.000"
Set-TextValue $ws "E19" "  -0.25%  "

# Row 20
Set-TextValue $ws "D20" "27.217.35"
Set-TextValue $ws "E20" "  -2.46%  "

# Row 21
Set-TextValue $ws "D21" "5.029"
Set-TextValue $ws "E21" "  -2.64%  "

# Row 22
Set-TextValue $ws "D22" "2.121.15"
Set-TextValue $ws "E22" "  -1.76%  "

# Row 23
Set-TextValue $ws "D23" "10.60"
Set-TextValue $ws "E23" "  -2.15%  "

# Row 24
Set-TextValue $ws "D24" "6.472"
Set-TextValue $ws "E24" "  -1.95%  "

# Row 25
Set-TextValue $ws "D25" "151.05"
Set-TextValue $ws "E25" "  -2.15%  "

# Row 26
Set-TextValue $ws "D26" "1.850"
Set-TextValue $ws "E26" "  -1.72%  "

# Row 27
Set-TextValue $ws "D27" "18.02"
Set-TextValue $ws "E27" "  -2.00%  "

# Row 28
Set-TextValue $ws "D28" "2.145"
Set-TextValue $ws "E28" "  -3.88%  "

# Row 29
Set-TextValue $ws "D29" "112.91"
Set-TextValue $ws "E29" "  -2.13%  "

# Row 30
Set-TextValue $ws "D30" "4.751"
Set-TextValue $ws "E30" "  -3.18%  "

# Row 31
Set-TextValue $ws "D31" "4.689"
Set-TextValue $ws "E31" "  +0.91%  "

# Row 32
Set-TextValue $ws "D32" "0.09008"
Set-TextValue $ws "E32" "  +0.37%  "

# Row 33
Set-TextValue $ws "D33" "0.05159"
Set-TextValue $ws "E33" "  -1.81%  "

# Row 34
Set-TextValue $ws "D34" "3.103"
Set-TextValue $ws "E34" "  -3.47%  "

# Row 35
Set-TextValue $ws "D35" "0.7525"
Set-TextValue $ws "E35" "  -1.79%  "

# Row 36
Set-TextValue $ws "D36" "1.172"
Set-TextValue $ws "E36" "  -4.92%  "

# Row 37
Set-TextValue $ws "D37" "0.02035"
Set-TextValue $ws "E37" "  -1.28%  "

# Row 38
Set-TextValue $ws "D38" "2.533"
Set-TextValue $ws "E38" "  -0.72%  "

# Row 39
Set-TextValue $ws "D39" "3.031"
Set-TextValue $ws "E39" "  +0.65%  "

# Row 40
Set-TextValue $ws "D40" "1.081"
Set-TextValue $ws "E40" "  -1.30%  "

# Row 41
Set-TextValue $ws "D41" "0.5345"
Set-TextValue $ws "E41" "  -4.38%  "

# Row 42
Set-TextValue $ws "D42" "6.643"
Set-TextValue $ws "E42" "  -4.56%  "

# Row 43
Set-TextValue $ws "D43" "115.03"
Set-TextValue $ws "E43" "  +3.28%  "

# Row 44
Set-TextValue $ws "D44" "8.486"
Set-TextValue $ws "E44" "  -0.05%  "

# Row 45
Set-TextValue $ws "D45" "0.1483"
Set-TextValue $ws "E45" "  -2.06%  "

# Row 46
Set-TextValue $ws "D46" "0.4663"
Set-TextValue $ws "E46" "  -3.08%  "

# Row 47
Set-TextValue $ws "D47" "1.001"
Set-TextValue $ws "E47" "  -0.16%  "

# Row 48
Set-TextValue $ws "D48" "10.12"
Set-TextValue $ws "E48" "  -4.22%  "

# Row 49
Set-TextValue $ws "D49" "1.572"
Set-TextValue $ws "E49" "  -3.65%  "

# Row 50
Set-TextValue $ws "D50" "64.95"
Set-TextValue $ws "E50" "  -3.85%  "

# Row 51
Set-TextValue $ws "D51" "36.41"
Set-TextValue $ws "E51" "  -1.48%  "
